# Add rural communities interactions back in, fix de_dg files.
# The "Legislature" row (row 12) has no interaction data recorded, so it is
# removed entirely (its row deleted, shifting everything below it up by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Legislature" row (row 12).
$ws.Rows(12).EntireRow.Delete()

# The two rows that now sit where the deletion happened (previously-row-11
# "Bureau of Reclamation" and the row that shifted up into row 12 "Friant
# Water Authority") are shrunk down to the default compact row height.
$ws.Rows(11).RowHeight = 15
$ws.Rows(12).RowHeight = 15

# Reflect the post-edit view/selection state: the row that was selected for
# deletion (now row 12) stays selected, and the viewport scrolls back up.
$ws.Range("A12:XFD12").Select()
$excel.ActiveWindow.ScrollRow = 11
